$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 35
$ws.Range("B1").Value = "clefairy"
$ws.Range("C1").Value = 113
$ws.Range("D1").Value = 6
$ws.Range("E1").Value = 56
$ws.Range("F1").Value = 75

$ws.Range("A2").Value = 13
$ws.Range("B2").Value = "weedle"
$ws.Range("C2").Value = 39
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 17
$ws.Range("F2").Value = 32
